$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '74.998.86'
$ws.Range("E2").Value = '  +1.05%  '

$ws.Range("D3").Value = '2.817.81'
$ws.Range("E3").Value = '  +6.75%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.58'
$ws.Range("E5").Value = '  +0.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.23'
$ws.Range("E6").Value = '  +1.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("E8").Value = '  +3.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.191'
$ws.Range("E9").Value = '  -4.28%  '

$ws.Range("D10").Value = '2.817.12'
$ws.Range("E10").Value = '  +6.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.376'
$ws.Range("E11").Value = '  +4.72%  '

$ws.Range("E12").Value = '  -2.04%  '

$ws.Range("E13").Value = '  +4.25%  '

$ws.Range("D14").Value = '3.332.53'
$ws.Range("E14").Value = '  +6.80%  '

$ws.Range("D15").Value = '74.986.43'
$ws.Range("E15").Value = '  +1.04%  '

$ws.Range("E16").Value = '  -1.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.90'
$ws.Range("E17").Value = '  +1.93%  '

$ws.Range("D18").Value = '2.813.73'
$ws.Range("E18").Value = '  +6.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.08'
$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.30'
$ws.Range("E20").Value = '  +3.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.34'
$ws.Range("E21").Value = '  +1.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.27'
$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("E23").Value = '  +0.42%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.80'
$ws.Range("E25").Value = '  +1.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.82'
$ws.Range("E26").Value = '  +5.61%  '

$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.18'
$ws.Range("E27").Value = '  +1.08%  '

$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.948.32'
$ws.Range("E28").Value = '  +6.09%  '

$ws.Range("E29").Value = '  +9.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '515.59'
$ws.Range("E31").Value = '  -1.54%  '

$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.65'
$ws.Range("E33").Value = '  -0.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.79'
$ws.Range("E34").Value = '  +1.86%  '

$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.04'
$ws.Range("E36").Value = '  +1.04%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.119'
$ws.Range("E37").Value = '  +1.82%  '

$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.95'
$ws.Range("E38").Value = '  +4.04%  '

$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '188.04'
$ws.Range("E39").Value = '  +16.32%  '

$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.37'
$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.342'
$ws.Range("E42").Value = '  +4.61%  '

$ws.Range("E43").Value = '  +1.84%  '

$ws.Range("E44").Value = '  +0.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.21'
$ws.Range("E45").Value = '  +2.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.07'
$ws.Range("E46").Value = '  +2.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0864'
$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("E48").Value = '  -2.19%  '

$ws.Range("E49").Value = '  +9.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.73'
$ws.Range("E50").Value = '  +3.21%  '

$ws.Range("E51").Value = '  +8.26%  '
